$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Price_id"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4160
$ws.Range("E1").Borders.LineStyle = 1

for ($r = 2; $r -le 427; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 4).Value2
}
